$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D, mirroring Code Review 2's structure
$ws.Range("D1").Value = "Code Review 3"

# Copy column C values (Code Review 2 scores) into new column D (Code Review 3)
$ws.Range("D2").Value = 25
$ws.Range("D3").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("D5").Value = 25
$ws.Range("D6").Value = 100

# Update column widths to match target layout
# (ColumnWidth is quantized internally to whole pixels, so these values
# are the closest settings that reproduce the target stored widths)
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 15.1

# Update the active selection as recorded in the edited workbook
$ws.Range("G8").Select()
